$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 410.2857
$ws.Range("I33").Value = 350.5
$ws.Range("J33").Value = 490
$ws.Range("K33").Value = 350.5
$ws.Range("L33").Value = 490
$ws.Range("M33").Value = -121.5
$ws.Range("N33").Value = -948
$ws.Range("H62").Value = 6001.3335
$ws.Range("I62").Value = 6005
$ws.Range("J62").Value = 5999.5
$ws.Range("K62").Value = 6005
$ws.Range("L62").Value = 5999.5
$ws.Range("M62").Value = -5381
$ws.Range("N62").Value = -7247.5
$ws.Range("H65").Value = 6001.3335
$ws.Range("I65").Value = 6005
$ws.Range("J65").Value = 5999.5
$ws.Range("K65").Value = 30025
$ws.Range("L65").Value = 29997.5
$ws.Range("M65").Value = -26905
$ws.Range("N65").Value = -36237.5
$ws.Range("H107").Value = 597.6818
$ws.Range("J107").Value = 674.2857
$ws.Range("L107").Value = 674.2857
$ws.Range("N107").Value = -4514.2857
$ws.Range("H112").Value = 55557824
$ws.Range("I112").Value = 333333920
$ws.Range("J112").Value = 2599.9333
$ws.Range("K112").Value = 1000001760
$ws.Range("L112").Value = 7799.7999
$ws.Range("M112").Value = -1000000652
$ws.Range("N112").Value = -10015.7999
$ws.Range("H115").Value = 415
$ws.Range("I115").Value = 415
$ws.Range("K115").Value = 1245
$ws.Range("M115").Value = 322
$ws.Range("H127").Value = 1209.1
$ws.Range("I127").Value = 884.4286
$ws.Range("K127").Value = 2653.2858
$ws.Range("M127").Value = 2306.7142
$ws.Range("H135").Value = 121766.47
$ws.Range("I135").Value = 94890.37
$ws.Range("J135").Value = 171039.33
$ws.Range("K135").Value = 854013.33
$ws.Range("L135").Value = 1539353.97
$ws.Range("M135").Value = -851478.33
$ws.Range("N135").Value = -1544423.97
$ws.Range("H138").Value = 5129821
$ws.Range("I138").Value = 1293.8695
$ws.Range("J138").Value = 12502079
$ws.Range("K138").Value = 3881.6085
$ws.Range("L138").Value = 37506237
$ws.Range("M138").Value = 1258.3915
$ws.Range("N138").Value = -37516517

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23120.486
$ws.Range("I32").Value = 23070.908
$ws.Range("J32").Value = 23264
$ws.Range("K32").Value = 23070.908
$ws.Range("L32").Value = 23264
$ws.Range("M32").Value = -22783.908
$ws.Range("N32").Value = -23838
$ws.Range("H61").Value = 55729.082
$ws.Range("I61").Value = 33536
$ws.Range("J61").Value = 170393.33
$ws.Range("K61").Value = 33536
$ws.Range("L61").Value = 170393.33
$ws.Range("M61").Value = -33324
$ws.Range("N61").Value = -170817.33
$ws.Range("H132").Value = 38585.21
$ws.Range("I132").Value = 20744.49
$ws.Range("J132").Value = 274974.75
$ws.Range("K132").Value = 62233.47
$ws.Range("L132").Value = 824924.25
$ws.Range("M132").Value = -59703.47
$ws.Range("N132").Value = -829984.25
$ws.Range("H136").Value = 55729.082
$ws.Range("I136").Value = 33536
$ws.Range("J136").Value = 170393.33
$ws.Range("K136").Value = 100608
$ws.Range("L136").Value = 511179.99
$ws.Range("M136").Value = -98058
$ws.Range("N136").Value = -516279.99

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2433.6562
$ws.Range("I134").Value = 2447.1936
$ws.Range("J134").Value = 2014
$ws.Range("K134").Value = 7341.5808
$ws.Range("L134").Value = 6042
$ws.Range("M134").Value = -4806.5808
$ws.Range("N134").Value = -11112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3106.1936
$ws.Range("I31").Value = 3121
$ws.Range("J31").Value = 3096.842
$ws.Range("K31").Value = 3121
$ws.Range("L31").Value = 3096.842
$ws.Range("M31").Value = -2826
$ws.Range("N31").Value = -3686.842
$ws.Range("H34").Value = 3106.1936
$ws.Range("I34").Value = 3121
$ws.Range("J34").Value = 3096.842
$ws.Range("K34").Value = 3121
$ws.Range("L34").Value = 3096.842
$ws.Range("M34").Value = -2919
$ws.Range("N34").Value = -3500.842
$ws.Range("H39").Value = 700.5
$ws.Range("I39").Value = 700.5
$ws.Range("K39").Value = 700.5
$ws.Range("M39").Value = -309.5
$ws.Range("H49").Value = 700.5
$ws.Range("I49").Value = 700.5
$ws.Range("K49").Value = 700.5
$ws.Range("M49").Value = -518.5
$ws.Range("H58").Value = 2382.9375
$ws.Range("I58").Value = 2394.3845
$ws.Range("J58").Value = 2333.3333
$ws.Range("K58").Value = 2394.3845
$ws.Range("L58").Value = 2333.3333
$ws.Range("M58").Value = -2191.3845
$ws.Range("N58").Value = -2739.3333
$ws.Range("H129").Value = 29330.8
$ws.Range("J129").Value = 29330.8
$ws.Range("L129").Value = 29330.8
$ws.Range("N129").Value = -39330.8
$ws.Range("H132").Value = 19308.229
$ws.Range("I132").Value = 1441.1892
$ws.Range("J132").Value = 52362.25
$ws.Range("K132").Value = 4323.5676
$ws.Range("L132").Value = 157086.75
$ws.Range("M132").Value = -1793.5676
$ws.Range("N132").Value = -162146.75
$ws.Range("H134").Value = 46709.207
$ws.Range("I134").Value = 1760.3846
$ws.Range("J134").Value = 99830.55
$ws.Range("K134").Value = 5281.1538
$ws.Range("L134").Value = 299491.65
$ws.Range("M134").Value = -2746.1538
$ws.Range("N134").Value = -304561.65
$ws.Range("H136").Value = 2382.9375
$ws.Range("I136").Value = 2394.3845
$ws.Range("J136").Value = 2333.3333
$ws.Range("K136").Value = 7183.1535
$ws.Range("L136").Value = 6999.999899999999
$ws.Range("M136").Value = -4633.1535
$ws.Range("N136").Value = -12099.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7987.9
$ws.Range("I3").Value = 5002.643
$ws.Range("J3").Value = 14953.5
$ws.Range("K3").Value = 15007.929
$ws.Range("L3").Value = 44860.5
$ws.Range("M3").Value = -14895.929
$ws.Range("N3").Value = -45084.5
$ws.Range("H103").Value = 2572.077
$ws.Range("J103").Value = 3520
$ws.Range("L103").Value = 10560
$ws.Range("N103").Value = -12318
$ws.Range("H131").Value = 1142.6666
$ws.Range("I131").Value = 463.16666
$ws.Range("J131").Value = 1414.4667
$ws.Range("K131").Value = 1389.49998
$ws.Range("L131").Value = 4243.4001
$ws.Range("M131").Value = 3650.50002
$ws.Range("N131").Value = -14323.4001
$ws.Range("H136").Value = 2885.5334
$ws.Range("I136").Value = 2218.75
$ws.Range("K136").Value = 6656.25
$ws.Range("M136").Value = -1556.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1047
$ws.Range("I102").Value = 996.93335
$ws.Range("J102").Value = 1197.2
$ws.Range("K102").Value = 996.93335
$ws.Range("L102").Value = 1197.2
$ws.Range("M102").Value = 625.06665
$ws.Range("N102").Value = -4441.2
$ws.Range("H107").Value = 625.7692
$ws.Range("I107").Value = 308.2
$ws.Range("J107").Value = 1684.3334
$ws.Range("K107").Value = 308.2
$ws.Range("L107").Value = 1684.3334
$ws.Range("M107").Value = 1611.8
$ws.Range("N107").Value = -5524.3334
$ws.Range("H122").Value = 1690.091
$ws.Range("I122").Value = 1337.4445
$ws.Range("J122").Value = 3277
$ws.Range("K122").Value = 4012.3335
$ws.Range("L122").Value = 9831
$ws.Range("M122").Value = -1562.3335
$ws.Range("N122").Value = -14731
$ws.Range("H128").Value = 55000
$ws.Range("J128").Value = 55000
$ws.Range("L128").Value = 55000
$ws.Range("N128").Value = -64960
$ws.Range("H132").Value = 54246
$ws.Range("I132").Value = 38477.406
$ws.Range("J132").Value = 92950.73
$ws.Range("K132").Value = 115432.218
$ws.Range("L132").Value = 278852.19
$ws.Range("M132").Value = -112902.218
$ws.Range("N132").Value = -283912.19

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9618872
$ws.Range("I7").Value = 17859878
$ws.Range("J7").Value = 4366.1665
$ws.Range("K7").Value = 17859878
$ws.Range("L7").Value = 4366.1665
$ws.Range("M7").Value = -17859766
$ws.Range("N7").Value = -4590.1665
$ws.Range("H82").Value = 901.1111
$ws.Range("I82").Value = 485
$ws.Range("J82").Value = 1733.3334
$ws.Range("K82").Value = 485
$ws.Range("L82").Value = 1733.3334
$ws.Range("M82").Value = -124
$ws.Range("N82").Value = -2455.3334
$ws.Range("H85").Value = 901.1111
$ws.Range("I85").Value = 485
$ws.Range("J85").Value = 1733.3334
$ws.Range("K85").Value = 485
$ws.Range("L85").Value = 1733.3334
$ws.Range("M85").Value = 763
$ws.Range("N85").Value = -4229.3334
$ws.Range("H126").Value = 9618872
$ws.Range("I126").Value = 17859878
$ws.Range("J126").Value = 4366.1665
$ws.Range("K126").Value = 53579634
$ws.Range("L126").Value = 13098.4995
$ws.Range("M126").Value = -53577164
$ws.Range("N126").Value = -18038.4995
$ws.Range("H132").Value = 47758.477
$ws.Range("I132").Value = 3920.3
$ws.Range("J132").Value = 81480.16
$ws.Range("K132").Value = 11760.9
$ws.Range("L132").Value = 244440.48
$ws.Range("M132").Value = -9230.900000000001
$ws.Range("N132").Value = -249500.48
$ws.Range("H136").Value = 49176.79
$ws.Range("I136").Value = 36904.723
$ws.Range("J136").Value = 74597.5
$ws.Range("K136").Value = 110714.169
$ws.Range("L136").Value = 223792.5
$ws.Range("M136").Value = -108164.169
$ws.Range("N136").Value = -228892.5
